# "fixed bad identifier for grammaticalGender"
#
# The table cell that used to read just "gender" must become
# "grammaticalGender" -- but split across two runs with the document's
# (single, unique) "_GoBack" bookmark sitting between "grammaticalG" and
# "ender" (this is where the author's cursor/last edit ended up). The
# "_GoBack" bookmark previously sat between the "http://www.tbxinfo.net/ns"
# and "/basic" runs near the top of the document; since a document can only
# carry one "_GoBack" bookmark, that old one is removed.

$d = $word.ActiveDocument

# --- 1. Drop the stale "_GoBack" bookmark (it used to live between the
#        namespace-URL run and the "/basic" run near the top of the file).
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# --- 2. Locate the table cell containing the lone word "gender".
$rng = $d.Content
$rng.Find.Execute("gender", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# --- 3. Replace it with two runs -- "grammaticalG" / "ender" -- keeping the
#        original run formatting (rFonts/color) on both, and re-creating the
#        original run's own rsidRPr on the first half (it's the same "logical"
#        run, just now shorter), with a fresh "_GoBack" bookmark pair sitting
#        right in between, exactly where the split happens.
$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00760C15" w:rsidRPr="002B2F0E" w:rsidRDefault="00760C15" w:rsidP="002B2F0E">
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="002B2F0E">
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>grammaticalG</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000"/>
              </w:rPr>
              <w:t>ender</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$rng.InsertXML($newXml)
